# Updated files from RMI Nov 2-24
$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")

# --- Text updates on the "About" sheet ---
# B6: "January 2020 and September 2020" -> "January 2020 and November 2020"
$wsAbout.Range("B6").Value = "January 2020 and November 2020"

# A27: "As of EPS 2.1.1, this variable is set up to model the impacts of the 2020"
#      -> "As of EPS 3.1, this variable is set up to model the impacts of the 2020"
$wsAbout.Range("A27").Value = "As of EPS 3.1, this variable is set up to model the impacts of the 2020"

# A28: "SARS-CoV-2 pandemic.  It uses the latest data available as of September 9,"
#      -> "SARS-CoV-2 pandemic.  It uses the latest data available as of November 10,"
$wsAbout.Range("A28").Value = "SARS-CoV-2 pandemic.  It uses the latest data available as of November 10,"

# --- "Data" sheet: update the source GDP figures (from the updated September -> November STEO) ---
$wsData.Range("A3").Value = "November STEO"
$wsData.Range("B3").Value = 19092
$wsData.Range("C3").Value = 18411
$wsData.Range("D3").Value = 19098

# Update the active selection on the Data sheet to B12, then restore the
# originally-active "About" sheet so the workbook's active tab is unchanged.
$wsData.Activate()
$wsData.Range("B12").Select()
$wsAbout.Activate()

$wb.Save()
